# Scoreboard.xlsx update
# - Fills in the per-lift detail columns (Minute/Second/Rep x2 + Snatch/CJ)
#   for ScoreM and ScoreF that were previously blank.
# - Replaces the generic "Team A".."Team F" placeholder names on the SFM /
#   SFF sheets with the real team names, and adds scores (Snatch / Clean &
#   Jerk totals) for each team, including a brand new 8th row on SFF.
# - Updates the active sheet / selection bookkeeping to match where the
#   author last left the cursor.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ScoreM ("ScoreM") - add lift detail columns C:Q for rows 2-11, fix a couple
# of QualifyPoints (column B) values that changed.
# ---------------------------------------------------------------------------
$wsScoreM = $wb.Worksheets.Item("ScoreM")

$wsScoreM.Cells.Item(2, 2).Value = 9
$wsScoreM.Cells.Item(4, 2).Value = 7

$scoreMRows = @(
    @{ Row = 2;  Vals = @(14,0,116,14,22,300,1,0,204,16,0,115,14,1,300) }
    @{ Row = 3;  Vals = @(14,0,107,13,39,300,1,0,218,16,0,116,13,36,300) }
    @{ Row = 4;  Vals = @(14,0,76,13,9,300,1,0,208,16,0,102,14,38,300) }
    @{ Row = 5;  Vals = @(14,0,119,14,36,300,1,0,190,16,0,110,14,14,300) }
    @{ Row = 6;  Vals = @(14,0,76,15,0,253,1,0,187,16,0,94,15,0,270) }
    @{ Row = 7;  Vals = @(14,0,84,13,10,300,1,0,185,16,0,118,15,0,259) }
    @{ Row = 8;  Vals = @(14,0,102,12,46,300,1,0,205,16,0,100,13,27,300) }
    @{ Row = 9;  Vals = @(14,0,113,12,27,300,1,0,224,16,0,116,14,23,300) }
    @{ Row = 10; Vals = @(14,0,107,14,50,300,1,0,217,16,0,75,15,0,290) }
    @{ Row = 11; Vals = @(14,0,104,14,48,300,1,0,210,16,0,110,13,49,300) }
)

foreach ($r in $scoreMRows) {
    $col = 3
    foreach ($v in $r.Vals) {
        $wsScoreM.Cells.Item($r.Row, $col).Value = $v
        $col++
    }
}

# ---------------------------------------------------------------------------
# ScoreF ("ScoreF") - add lift detail columns C:N for rows 2-10.
# ---------------------------------------------------------------------------
$wsScoreF = $wb.Worksheets.Item("ScoreF")

$scoreFRows = @(
    @{ Row = 2;  Vals = @(16,0,82,13,28,300,1,0,181,11,0,100) }
    @{ Row = 3;  Vals = @(16,0,88,15,0,293,1,0,183,11,0,76) }
    @{ Row = 4;  Vals = @(16,0,82,13,12,300,1,0,215,11,0,92) }
    @{ Row = 5;  Vals = @(16,0,100,12,58,300,1,0,214,11,0,83) }
    @{ Row = 6;  Vals = @(16,0,106,14,11,300,1,0,182,11,0,75) }
    @{ Row = 7;  Vals = @(16,0,84,14,12,300,1,0,223,11,0,93) }
    @{ Row = 8;  Vals = @(16,0,82,13,57,300,1,0,188,11,0,81) }
    @{ Row = 9;  Vals = @(16,0,94,12,20,300,1,0,186,11,0,85) }
    @{ Row = 10; Vals = @(16,0,112,14,35,300,1,0,189,11,0,104) }
)

foreach ($r in $scoreFRows) {
    $col = 3
    foreach ($v in $r.Vals) {
        $wsScoreF.Cells.Item($r.Row, $col).Value = $v
        $col++
    }
}

# ---------------------------------------------------------------------------
# SFM - replace "Team A".."Team F" placeholders with real team names and
# fill in the Snatch / Clean & Jerk score columns.
# ---------------------------------------------------------------------------
$wsSFM = $wb.Worksheets.Item("SFM")

$sfmRows = @(
    @{ Row = 2; Name = "Magnus Øslebye og Vegard Tangen";              B = 75; C = 100 }
    @{ Row = 3; Name = "Anders J. Svalestuen og Gabriel Kristiansen";  B = 78; C = 111 }
    @{ Row = 4; Name = "Kasper Støen Nerbøvik og Håvard Idland";       B = 80; C = 99 }
    @{ Row = 5; Name = "Magnus Ødegaard og Kornelius Skrettingland";   B = 66; C = 85 }
    @{ Row = 6; Name = "Ole Andre Elvebakk og Georg Kongsvik";         B = 52; C = 98 }
    @{ Row = 7; Name = "Håkon Konningen og Njål Christensen";          B = 77; C = 104 }
)

foreach ($r in $sfmRows) {
    $wsSFM.Cells.Item($r.Row, 1).Value = $r.Name
    $wsSFM.Cells.Item($r.Row, 2).Value = $r.B
    $wsSFM.Cells.Item($r.Row, 3).Value = $r.C
}

# ---------------------------------------------------------------------------
# SFF - replace "Team A".."Team F" placeholders with real team names, fill
# in the score columns, and add a brand new 8th team row.
# ---------------------------------------------------------------------------
$wsSFF = $wb.Worksheets.Item("SFF")

$sffRows = @(
    @{ Row = 2; Name = "Renate Berntsen Hansen og Karoline Granås";     B = 66; C = 90 }
    @{ Row = 3; Name = "Maria Hanssen og Cecilie Rabben";               B = 59; C = 85 }
    @{ Row = 4; Name = "Victoria Christensen og Helene Rye Martinsen"; B = 72; C = 81 }
    @{ Row = 5; Name = "Marianne U. Henriksen og Mari S. Andersen";    B = 78; C = 75 }
    @{ Row = 6; Name = "Dawn Stewart og Marie Vik";                    B = 55; C = 91 }
    @{ Row = 7; Name = "Sara Yuzer og Martine Baalsrud";                B = 69; C = 80 }
    @{ Row = 8; Name = "Frid Kaspersen og Renate Loraas";               B = 70; C = 78 }
)

foreach ($r in $sffRows) {
    $wsSFF.Cells.Item($r.Row, 1).Value = $r.Name
    $wsSFF.Cells.Item($r.Row, 2).Value = $r.B
    $wsSFF.Cells.Item($r.Row, 3).Value = $r.C
}

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping - match where the author's cursor
# ended up. ScoreF must be activated LAST so it becomes the workbook's
# active tab.
# ---------------------------------------------------------------------------
$wsScoreM.Activate()
$wsScoreM.Range("S13").Select()

$wsSFM.Activate()
$wsSFM.Range("B17").Select()

$wsSFF.Activate()
$wsSFF.Range("J10").Select()

$wsScoreF.Activate()
$wsScoreF.Range("J20").Select()
